# Apply the edits described by the diff:
# - Companies sheet: A2 "test123" -> "jkugjyh"; delete row 3 (sgrhdtjgf / TRUE)
# - Locations sheet: delete row 2 (BC / test123)
# - AssetTypes sheet: delete row 2 (BC / BC / #562c33)

$wb = $excel.ActiveWorkbook

$wsCompanies = $wb.Worksheets.Item("Companies")
$wsLocations = $wb.Worksheets.Item("Locations")
$wsAssetTypes = $wb.Worksheets.Item("AssetTypes")

# Companies: rename test123 -> jkugjyh, then drop the now-redundant sgrhdtjgf row
$wsCompanies.Range("A2").Value = "jkugjyh"
$wsCompanies.Rows.Item(3).Delete()

# Locations: drop the BC / test123 row
$wsLocations.Rows.Item(2).Delete()

# AssetTypes: drop the BC / BC / #562c33 row
$wsAssetTypes.Rows.Item(2).Delete()
